$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-19 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-20 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("69×45=3105", $true, $false, $false, $false, $false, $true, 1, $false, "91×97=8827", 2) | Out-Null
$d.Content.Find.Execute("20×34=680", $true, $false, $false, $false, $false, $true, 1, $false, "86×43=3698", 2) | Out-Null
$d.Content.Find.Execute("18×41=738", $true, $false, $false, $false, $false, $true, 1, $false, "22×47=1034", 2) | Out-Null
$d.Content.Find.Execute("25×28=700", $true, $false, $false, $false, $false, $true, 1, $false, "65×82=5330", 2) | Out-Null
$d.Content.Find.Execute("37×82=3034", $true, $false, $false, $false, $false, $true, 1, $false, "82×15=1230", 2) | Out-Null
$d.Content.Find.Execute("76×17=1292", $true, $false, $false, $false, $false, $true, 1, $false, "21×84=1764", 2) | Out-Null
$d.Content.Find.Execute("42×11=462", $true, $false, $false, $false, $false, $true, 1, $false, "80×23=1840", 2) | Out-Null
$d.Content.Find.Execute("20×80=1600", $true, $false, $false, $false, $false, $true, 1, $false, "56×44=2464", 2) | Out-Null
$d.Content.Find.Execute("13×11=143", $true, $false, $false, $false, $false, $true, 1, $false, "92×98=9016", 2) | Out-Null
$d.Content.Find.Execute("82×84=6888", $true, $false, $false, $false, $false, $true, 1, $false, "25×36=900", 2) | Out-Null
$d.Content.Find.Execute("31×94=2914", $true, $false, $false, $false, $false, $true, 1, $false, "98×41=4018", 2) | Out-Null
$d.Content.Find.Execute("77×70=5390", $true, $false, $false, $false, $false, $true, 1, $false, "59×80=4720", 2) | Out-Null
$d.Content.Find.Execute("59×12=708", $true, $false, $false, $false, $false, $true, 1, $false, "19×35=665", 2) | Out-Null
$d.Content.Find.Execute("87×98=8526", $true, $false, $false, $false, $false, $true, 1, $false, "52×94=4888", 2) | Out-Null
$d.Content.Find.Execute("26×35=910", $true, $false, $false, $false, $false, $true, 1, $false, "88×65=5720", 2) | Out-Null
$d.Content.Find.Execute("35×63=2205", $true, $false, $false, $false, $false, $true, 1, $false, "67×91=6097", 2) | Out-Null
$d.Content.Find.Execute("69×35=2415", $true, $false, $false, $false, $false, $true, 1, $false, "15×23=345", 2) | Out-Null
$d.Content.Find.Execute("47×14=658", $true, $false, $false, $false, $false, $true, 1, $false, "43×98=4214", 2) | Out-Null
$d.Content.Find.Execute("14×40=560", $true, $false, $false, $false, $false, $true, 1, $false, "94×45=4230", 2) | Out-Null
$d.Content.Find.Execute("75×33=2475", $true, $false, $false, $false, $false, $true, 1, $false, "94×94=8836", 2) | Out-Null
$d.Content.Find.Execute("57×66=3762", $true, $false, $false, $false, $false, $true, 1, $false, "35×65=2275", 2) | Out-Null
$d.Content.Find.Execute("41×49=2009", $true, $false, $false, $false, $false, $true, 1, $false, "36×30=1080", 2) | Out-Null
$d.Content.Find.Execute("12×97=1164", $true, $false, $false, $false, $false, $true, 1, $false, "62×17=1054", 2) | Out-Null
$d.Content.Find.Execute("44×19=836", $true, $false, $false, $false, $false, $true, 1, $false, "83×22=1826", 2) | Out-Null
$d.Content.Find.Execute("67×74=4958", $true, $false, $false, $false, $false, $true, 1, $false, "88×78=6864", 2) | Out-Null
